$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4+ down by one.
$ws.Rows("4:4").Insert()

# Update row 2: "when" becomes the numeric year 2019 (was text "2017 - Present").
$ws.Range("B2").Value = 2019

# Update row 3: "when" becomes "2017-2018" and "why" becomes
# "Quantitative Methods II (Psychology MSc)." (same text as row 2).
$ws.Range("B3").Value = "2017-2018"
$ws.Range("E3").Value = "Quantitative Methods II (Psychology MSc)."

# Fill in the newly inserted row 4 with the split-off teaching entry.
$ws.Range("B4").Value = 2017
$ws.Range("E4").Value = "Quantitative Methods I (Psychology MSc)."

# Match the author's final selection position.
$ws.Range("C19").Select()
